# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# 1. The reporting period label changes from 2507 to 2508 (shared by the
#    existing worker row).
# 2. A new worker (MANUEL ALFONSO FUENTES QUINTANA, CC 73212454) is added
#    as a new data row right below the existing one, duplicating the row's
#    formatting.
# 3. The totals (Valor Mora / Cant. Trabajadores) are updated to reflect
#    the second worker.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- 1. Update the period text (also updates row 16, which references the
#        same text) ---------------------------------------------------
$ws.Range("E16").Value2 = "2508"

# --- 2. Insert a new row for the second worker, copying row 16's format --
$ws.Rows("17:17").Insert(-4121)            # xlShiftDown
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B17").Value2 = "CC"
$ws.Range("C17").Value2 = "73212454"
$ws.Range("D17").Value2 = "MANUEL ALFONSO FUENTES QUINTANA"
$ws.Range("E17").Value2 = "2508"
$ws.Range("F17").Value2 = 56940
$ws.Range("G17").Value2 = 1423500

# --- 3. Update the summary figures -------------------------------------
$ws.Range("E11").Value2 = 113880
$ws.Range("C13").Value2 = 2

Write-Host "Applied EC update: added MANUEL ALFONSO FUENTES QUINTANA, period 2508"
